$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize Spanish connector words (de/del/la/los/el) to Title Case, plus two special-case renames ---
$ws.Range("B4").Value = "Rincón De Romos"
$ws.Range("B5").Value = "San Francisco De Los Romo"
$ws.Range("B17").Value = "Amatenango De La Frontera"
$ws.Range("B25").Value = "Marqués De Comillas"
$ws.Range("B26").Value = "Mazapa De Madero"
$ws.Range("B44").Value = "Hidalgo Del Parral"
$ws.Range("B51").Value = "Valle De Zaragoza"
$ws.Range("B58").Value = "San Juan De Sabinas"
$ws.Range("A67").Value = "Ciudad De México"
$ws.Range("B92").Value = "Pánuco De Coronado"
$ws.Range("B95").Value = "San Pedro Del Gallo"
$ws.Range("A100").Value = "Estado De México"
$ws.Range("B101").Value = "Almoloya De Juárez"
$ws.Range("B105").Value = "Atizapán De Zaragoza"
$ws.Range("B109").Value = "Ecatepec De Morelos"
$ws.Range("B110").Value = "Ixtapan De La Sal"
$ws.Range("B112").Value = "Naucalpan De Juárez"
$ws.Range("B115").Value = "San Felipe Del Progreso"
$ws.Range("B116").Value = "Soyaniquilpan De Juárez"
$ws.Range("B123").Value = "Tlalnepantla De Baz"
$ws.Range("B126").Value = "Villa De Allende"
$ws.Range("B133").Value = "Apaseo El Alto"
$ws.Range("B134").Value = "Apaseo El Grande"
$ws.Range("B141").Value = "Jaral Del Progreso"
$ws.Range("B147").Value = "Purísima Del Rincón"
$ws.Range("B152").Value = "San Francisco Del Rincón"
$ws.Range("B153").Value = "San Luis De La Paz"
$ws.Range("B154").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B155").Value = "Silao De La Victoria"
$ws.Range("B157").Value = "Valle De Santiago"
$ws.Range("B160").Value = "Acapulco De Juárez"
$ws.Range("B161").Value = "Ajuchitlán Del Progreso"
$ws.Range("B164").Value = "Atoyac De Álvarez"
$ws.Range("B165").Value = "Ayutla De Los Libres"
$ws.Range("B168").Value = "Buenavista De Cuéllar"
$ws.Range("B169").Value = "Chilapa De Álvarez"
$ws.Range("B170").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B173").Value = "Coyuca De Benítez"
$ws.Range("B174").Value = "Cutzamala De Pinzón"
$ws.Range("B179").Value = "Iguala De La Independencia"
$ws.Range("B181").Value = "Zihuatanejo De Azueta"
$ws.Range("B183").Value = "Mártir De Cuilapan"
$ws.Range("B192").Value = "Técpan De Galeana"
$ws.Range("B196").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B197").Value = "Tlapa De Comonfort"
$ws.Range("B203").Value = "Agua Blanca De Iturbide"
$ws.Range("B204").Value = "Atotonilco El Grande"
$ws.Range("B206").Value = "Cuautepec De Hinojosa"
$ws.Range("B207").Value = "Huasca De Ocampo"
$ws.Range("B208").Value = "Huejutla De Reyes"
$ws.Range("B213").Value = "Mineral Del Chico"
$ws.Range("B214").Value = "Pachuca De Soto"
$ws.Range("B217").Value = "Tulancingo De Bravo"
$ws.Range("B221").Value = "Acatlán De Juárez"
$ws.Range("B224").Value = "Atemajac De Brizuela"
$ws.Range("B226").Value = "Autlán De Navarro"
$ws.Range("B236").Value = "Ixtlahuacán Del Río"
$ws.Range("B240").Value = "Jilotlán De Los Dolores"
$ws.Range("B243").Value = "La Manzanilla De La Paz"
$ws.Range("B244").Value = "Lagos De Moreno"
$ws.Range("B252").Value = "San Cristóbal De La Barranca"
$ws.Range("B253").Value = "San Diego De Alejandría"
$ws.Range("B255").Value = "San Miguel El Alto"
$ws.Range("B256").Value = "San Sebastián Del Oeste"
$ws.Range("B257").Value = "Santa María De Los Ángeles"
$ws.Range("B258").Value = "Tamazula De Gordiano"
$ws.Range("B263").Value = "Tepatitlán De Morelos"
$ws.Range("B265").Value = "Tizapán El Alto"
$ws.Range("B269").Value = "Unión De San Antonio"
$ws.Range("B270").Value = "Unión De Tula"
$ws.Range("B271").Value = "Valle De Juárez"
$ws.Range("B272").Value = "Yahualica De González Gallo"
$ws.Range("B273").Value = "Zacoalco De Torres"
$ws.Range("B275").Value = "Zapotlán El Grande"
$ws.Range("B286").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B288").Value = "Cojumatlán De Régules"
$ws.Range("B343").Value = "Coatlán Del Río"
$ws.Range("B353").Value = "Tlaltizapán De Zapata"
$ws.Range("B359").Value = "Amatlán De Cañas"
$ws.Range("B373").Value = "San Nicolás De Los Garza"
$ws.Range("B375").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B378").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B380").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B381").Value = "Oaxaca De Juárez"
$ws.Range("B385").Value = "San Antonino El Alto"
$ws.Range("B386").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B395").Value = "San Miguel Del Puerto"
$ws.Range("B413").Value = "Zimatlán De Álvarez"
$ws.Range("B421").Value = "Izúcar De Matamoros"
$ws.Range("B422").Value = "Palmar De Bravo"
$ws.Range("B425").Value = "Tecali De Herrera"
$ws.Range("B436").Value = "San Juan Del Río"
$ws.Range("B448").Value = "Santa María Del Río"
$ws.Range("B453").Value = "Villa De Ramos"
$ws.Range("B464").Value = "Nacozari De García"
$ws.Range("B492").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B502").Value = "Cosamaloapan De Carpio"
$ws.Range("B507").Value = "Ignacio De La Llave"
$ws.Range("B511").Value = "Juchique De Ferrer"
$ws.Range("B513").Value = "Martínez De La Torre"
$ws.Range("B514").Value = "Medellín De Bravo"
$ws.Range("B521").Value = "Paso De Ovejas"
$ws.Range("B523").Value = "Poza Rica De Hidalgo"
$ws.Range("B528").Value = "Soledad De Doblado"
$ws.Range("B534").Value = "Vega De Alatorre"
$ws.Range("B545").Value = "Concepción Del Oro"
$ws.Range("B546").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B556").Value = "Mezquital Del Oro"
$ws.Range("B566").Value = "Villa De Cos"
$ws.Range("A130").Value = "Guanajuato"
$ws.Range("B370").Value = "Montemorelos"

# --- Fix tiny floating point precision differences ---
$ws.Range("D303").Value = 0.009703196347031965
$ws.Range("D547").Value = 0.009703196347031965

# --- Remove trailing footnote rows (574-578); shrinks used range to A1:D572 ---
$ws.Range("A574:A578").EntireRow.Delete() | Out-Null

